$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.898.11'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '1.551.64'
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.67'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("E6").Value = '  +0.42%  '
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.71'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +1.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.248'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0858'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '1.772.92'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("D13").Value = '1.547.31'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.73'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("E15").Value = '  +1.93%  '
$ws.Range("D16").Value = '26.886.06'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.65'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +0.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '216.88'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +2.29%  '
$ws.Range("D19").Value = '0.0₃0689'
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.22'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("E21").Value = '  +0.56%  '
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.21'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.83'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.60'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.88'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("E29").Value = '  +1.11%  '
$ws.Range("E30").Value = '  +2.76%  '
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("D33").Value = '1.427.16'
$ws.Range("E33").Value = '  +5.17%  '
$ws.Range("E34").Value = '  +3.17%  '
$ws.Range("E35").Value = '  +4.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.959'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +2.64%  '
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.522'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.70'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.27'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +3.99%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.985'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.54'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("D47").Value = '1.686.94'
$ws.Range("E47").Value = '  +1.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.16'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +1.04%  '
$ws.Range("E49").Value = '  +3.43%  '
$ws.Range("E50").Value = '  +4.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0956'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +1.49%  '
